$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update the Russian title - drop the period after "5.6.1.1"
$ws.Range("B1").Value = "5.6.1.1 Доля замужних женщин и сексуально активных не замужних женщин в возрасте 15-49 лет, которые были осведомлены о соврменном методе контрацепции"

# Row 6: "urban" -> new wording
$ws.Range("A6").Value = "Шаар жерлери"
$ws.Range("B6").Value = "Городские поселения"
$ws.Range("C6").Value = "City"

# Row 7: "rural" -> new wording
$ws.Range("A7").Value = "Айыл аймагы"
$ws.Range("B7").Value = "Сельская местность"
$ws.Range("C7").Value = "Village"

# Update the selected range shown in the sheet view
$ws.Range("A6:C7").Select()
